# Applies the "alteracao na documentacao da API" edit:
#  1. Removes the stray empty paragraph that sat just above the "URL" heading.
#  2. Reflows the "A URL utilizada e ..." paragraph: drop its leading tab run
#     and give the paragraph a first-line indent instead.
#  3. Adds a new "Bibliotecas utilizadas" heading (same style as the other
#     3rd-level numbered headings) followed by its explanatory paragraph.

$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, [string]$text) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# --- Step 1: delete the empty paragraph directly above the "URL" heading ---
$urlHeadingIdx = Get-ParaIndexByText $d "URL"
$emptyBeforeUrlIdx = $urlHeadingIdx - 1
$emptyBeforeUrl = $d.Paragraphs.Item($emptyBeforeUrlIdx)
$emptyBeforeUrl.Range.Delete()

# --- Step 2: fix up the "A URL utilizada e ..." paragraph ---
$urlHeadingIdx = Get-ParaIndexByText $d "URL"
$urlBodyIdx = $urlHeadingIdx + 1
$urlBody = $d.Paragraphs.Item($urlBodyIdx)

# Remove the leading tab character (first character of the paragraph range)
# and replace it with a first-line indent on the paragraph itself.
$firstChar = $urlBody.Range.Characters.Item(1)
if ($firstChar.Text -eq [char]9) {
    $firstChar.Delete()
}
$urlBody.Format.FirstLineIndent = 36

# --- Step 3: insert the two new paragraphs after the "A URL utilizada..." one ---
$urlBody = $d.Paragraphs.Item($urlBodyIdx)
$urlBody.Range.InsertParagraphAfter()
$newHeadingPara = $d.Paragraphs.Item($urlBodyIdx + 1)
$newHeadingPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($urlBodyIdx + 1)
$headingRange = $d.Range($headingPara.Range.Start, $headingPara.Range.End)
$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:shd w:val="clear" w:fill="ffffff"/><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:left="2160" w:hanging="360"/><w:jc w:val="both"/><w:rPr><w:b w:val="1"/><w:color w:val="3b4151"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:u w:val="none"/></w:rPr></w:pPr><w:r><w:rPr><w:b w:val="1"/><w:color w:val="3b4151"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">Bibliotecas utilizadas</w:t></w:r></w:p>'
$headingRange.InsertXML($headingXml)

$bodyPara = $d.Paragraphs.Item($urlBodyIdx + 2)
$bodyRange = $d.Range($bodyPara.Range.Start, $bodyPara.Range.End)
$bodyXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:shd w:val="clear" w:fill="ffffff"/><w:spacing w:after="200" w:line="240" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:color w:val="3b4151"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="3b4151"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:rtl w:val="0"/></w:rPr><w:tab/><w:t xml:space="preserve">Para a cria&#231;&#227;o desta API, devem ser importadas as bibliotecas Express (para desenvolvimento da API e suas rotas), Sequelize (para fazer a conex&#227;o com o banco de dados MySQL) e Jest (para realizar os testes no c&#243;digo).</w:t></w:r></w:p>'
$bodyRange.InsertXML($bodyXml)

Write-Output "edit applied"
